$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order corresponds to: C, D, E, F, G, H, J, L
$colIndexes = @(3, 4, 5, 6, 7, 8, 10, 12)

$data = @(
    @(0.3740083992085772, 0.02724418757605918, 0.1401897643330425, 6.557428756353914, 6.01939653040273, 3.997281369321001, 0.3938509743047121, 0.1118087815775866),
    @(0.372667010370364, 0.02489096839891403, 0.1406261917705685, 6.434948212661141, 5.86786478500062, 3.938905539795996, 0.3915715013672951, 0.112451925594943),
    @(0.3720460413837259, 0.02343810531269952, 0.1409525226156347, 6.364209360459199, 5.779160603376397, 3.905675968507353, 0.3903992367865072, 0.112876709259158),
    @(0.3718436792803885, 0.02284387653210729, 0.1411001577118274, 6.336493176804623, 5.744089390514773, 3.89278550608401, 0.3899783203635536, 0.1130573438889702),
    @(0.3718131309068298, 0.0227450690194857, 0.1411255567447345, 6.331957644482941, 5.738330452565265, 3.890684183614781, 0.3899118467763145, 0.1130877934986927),
    @(0.3720431073732158, 0.02343010036373983, 0.1409544543752137, 6.363831089934934, 5.778683280535063, 3.905499496124378, 0.3903933307029206, 0.1128791148435759),
    @(0.3735036620075931, 0.02643434981811765, 0.1403281167795853, 6.514263321210535, 5.966239611626577, 3.976607092242318, 0.3930175961314717, 0.1120243431399821),
    @(0.3779889155194525, 0.03227036071579903, 0.1395641598401305, 6.84527292061199, 6.369124433330228, 4.137089045987068, 0.3999858260432489, 0.110584625320282),
    @(0.3822922488738811, 0.03653591584325966, 0.1392878426893844, 7.111286806379781, 6.68750992840728, 4.268271190672522, 0.4062432343923064, 0.1096701222097405),
    @(0.3844733748428268, 0.03847420542565771, 0.1392244798864475, 7.237455936493916, 6.837438603795249, 4.330934473592777, 0.4093431699045738, 0.1092850110907602),
    @(0.3853318263374206, 0.0392080896274507, 0.1392094891314279, 7.285990733471579, 6.894963572998506, 4.355101455187423, 0.4105539854504485, 0.1091436088605633),
    @(0.3851454917516719, 0.03905003587208, 0.1392123165749446, 7.275503962648486, 6.882540885201763, 4.349877061900884, 0.4102915637442379, 0.1091738654864862),
    @(0.3845433465443762, 0.03853458342400273, 0.1392230659115832, 7.241433655332344, 6.842156077834773, 4.332913885022492, 0.4094420412343425, 0.1092732891105967),
    @(0.3841787595913218, 0.03821884648347407, 0.1392308239118307, 7.220663700067178, 6.817517457855161, 4.322580703489507, 0.4089265095693975, 0.1093347656317682),
    @(0.3821542357732142, 0.03640921834467292, 0.1392932412752259, 7.103146468933744, 6.677815739129471, 4.264236807723535, 0.4060457876422845, 0.109695910854068),
    @(0.3809697593354144, 0.03529866632503342, 0.1393475253000673, 7.0323854684423, 6.593430308114989, 4.229215660517411, 0.4043438101736712, 0.1099253671960554),
    @(0.3803094957085307, 0.03465969297042903, 0.1393846128102858, 6.992170153445642, 6.545372367118489, 4.209353153737027, 0.4033886989461735, 0.1100602536845887),
    @(0.3800895404093154, 0.03444330486855307, 0.1393981761159218, 6.978636693420754, 6.529182394438465, 4.202676023481388, 0.4030693911606562, 0.1101064239887997),
    @(0.3810936712133639, 0.03541690709790402, 0.1393411394118047, 7.039867838568057, 6.602363649254983, 4.232914609046873, 0.404522519634483, 0.1099006401705171),
    @(0.3847193259941548, 0.03868598548218927, 0.1392196638944583, 7.251420251306001, 6.853997568860393, 4.337884434673185, 0.4096905600657408, 0.1092439658241133),
    @(0.3872785531148395, 0.04082200908493405, 0.1391927712490926, 7.394103624924924, 7.022836498142397, 4.409043827508299, 0.4132837161555187, 0.1088406128405772),
    @(0.3858951678995481, 0.03968195115021445, 0.1392023070030532, 7.317540962805481, 6.932317033349364, 4.370828136050818, 0.4113460896759733, 0.1090535312789811),
    @(0.3810375861814634, 0.03536345202419966, 0.1393440081633379, 7.036483608324033, 6.598323465581927, 4.231241468999315, 0.4044416523065877, 0.1099118100100327),
    @(0.3765999893390841, 0.03069630908394316, 0.1397209344320594, 6.75178270529554, 6.256275961390372, 4.091376358756747, 0.3979029319190204, 0.1109488849938707)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $colIndexes.Count; $j++) {
        $c = $colIndexes[$j]
        $ws.Cells.Item($r, $c).Value2 = $rowValues[$j]
    }
}

Write-Host "Updated $($data.Count) rows"
